# Update workbook to reflect data refresh through 2022-05-06
# (adds one more day of carjacking data: 2022-05-14 commit / new daily totals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet tab and update the "Through" date reference
$ws.Name = "Through 2022-05-06"

# Update the header cell text for the current (partial) month column
$ws.Range("B1").Value = "May 2022 (through May 06)"

# Updated / new cell values (monthly counts by neighborhood), all within
# the "May" column of each year (columns B, G, L, Q, V, AA, AF, AK).
$updates = @{
    "B2"   = 2
    "AA3"  = 2
    "AF3"  = 1
    "L4"   = 3
    "L5"   = 1
    "V5"   = 2
    "Q7"   = 1
    "AF7"  = 1
    "Q11"  = 1
    "L20"  = 1
    "G22"  = 1
    "L25"  = 1
    "G26"  = 1
    "Q31"  = 1
    "B46"  = 1
    "Q68"  = 1
    "V70"  = 1
    "G77"  = 1
    "AK94" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
